# Generate Report for Handback
# Updates the localization-status workbook to reflect that the zh-cn and
# de-de handback packages came back in sync with en-US:
#  - "Status" cells flip from "Ready for handoff" to
#    "Handed back: in sync with en-US"
#  - The per-language "Latest Target File" / "Latest Handback File" /
#    "Latest Handback DateTime" columns are populated
#  - A "Latest Target File" hyperlink (a.md) is added for each data row
#  - A couple of column widths are widened to fit the new content

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status columns (E, F) for both rows ---
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# Widen the zh-cn / de-de status columns to fit the longer text
$overview.Columns.Item(5).ColumnWidth = 29.1666666666667
$overview.Columns.Item(6).ColumnWidth = 29.1666666666667

# --- zh-cn sheet -----------------------------------------------------
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$zhcnTarget = "https://github.com/OpenLocalizationTestOrg/oltest/blob/c7f9934434264f51fb883a95351ab07c1353f5c6/e2e/a.md"
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $zhcnTarget, "", "", "a.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $zhcnTarget, "", "", "a.md") | Out-Null

$zhcn.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

$zhcn.Range("K2").Value = "2016-08-12 02:59:06"
$zhcn.Range("K3").Value = "2016-08-12 02:59:06"

# Widen Status (C) and Latest Handback File (J) columns
$zhcn.Columns.Item(3).ColumnWidth = 29.1666666666667
$zhcn.Columns.Item(10).ColumnWidth = 39.1666666666667

# --- de-de sheet -------------------------------------------------------
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

$dedeTarget = "https://github.com/OpenLocalizationTestOrg/oltest/blob/c7f9934434264f51fb883a95351ab07c1353f5c6/e2e/a.md"
$dede.Hyperlinks.Add($dede.Range("I2"), $dedeTarget, "", "", "a.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("I3"), $dedeTarget, "", "", "a.md") | Out-Null

$dede.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$dede.Range("K2").Value = "2016-08-12 02:59:13"
$dede.Range("K3").Value = "2016-08-12 02:59:13"

# Widen Status (C) and Latest Handback File (J) columns
$dede.Columns.Item(3).ColumnWidth = 29.1666666666667
$dede.Columns.Item(10).ColumnWidth = 39.1666666666667
